# Updated cryptos list — applies the per-cell value edits described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '37.345.22'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.05%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.036.71'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.54%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.65'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '59.09'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.01%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.393'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +4.05%  '
$ws.Range('E10').Value = '  +2.42%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.14'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +6.64%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.337.22'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.64%  '
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '22.03'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.035.49'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.80%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '37.283.70'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.38'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.27'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.52%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '229.88'
$ws.Range('D22').Style = "Normal"
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +4.52%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.35'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.19%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '164.66'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.139'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.82%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.85'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.80%  '
$ws.Range('E30').Value = '  +3.21%  '
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0678'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +10.55%  '
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.54'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +11.40%  '
$ws.Range('E35').Value = '  +0.54%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.59'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +5.84%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +2.16%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.45'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  +3.00%  '
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0216'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.42%  '
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.64'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.25%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.395.81'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.90%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '91.57'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('E47').Value = '  +3.41%  '
$ws.Range('E48').Value = '  +3.32%  '
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.228.75'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.66%  '
